# Apply cryptos list update (commit: Mon Jun 19 23:53:00 UTC 2023, GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.731.94"
$ws.Range("E2").Value = "  +1.40%  "

# Row 3
$ws.Range("D3").Value = "1.731.36"
$ws.Range("E3").Value = "  +0.61%  "

# Row 4
$ws.Range("E4").Value = "  -0.26%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.18"
$ws.Range("E5").Value = "  -0.83%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9984"
$ws.Range("E6").Value = "  -0.24%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4919"
$ws.Range("E7").Value = "  +0.95%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2622"
$ws.Range("E8").Value = "  +0.48%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06218"
$ws.Range("E9").Value = "  +0.01%  "

# Row 10
$ws.Range("D10").Value = "1.729.77"
$ws.Range("E10").Value = "  +0.32%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "15.96"
$ws.Range("E11").Value = "  +3.34%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06980"
$ws.Range("E12").Value = "  -0.51%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6123"
$ws.Range("E13").Value = "  +2.72%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.511"
$ws.Range("E14").Value = "  -0.45%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.22"
$ws.Range("E15").Value = "  +0.04%  "

# Row 16
$ws.Range("E16").Value = "  -0.22%  "

# Row 17
$ws.Range("D17").Value = "26.530.10"
$ws.Range("E17").Value = "  +0.54%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9981"
$ws.Range("E18").Value = "  -0.28%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007201"
$ws.Range("E19").Value = "  -0.33%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.42"
$ws.Range("E20").Value = "  +0.68%  "

# Row 21
$ws.Range("D21").Value = "1.951.42"
$ws.Range("E21").Value = "  +0.14%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.464"
$ws.Range("E22").Value = "  -0.27%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.568"
$ws.Range("E23").Value = "  +0.09%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.109"
$ws.Range("E24").Value = "  -1.05%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.05"
$ws.Range("E25").Value = "  +0.52%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.32"
$ws.Range("E26").Value = "  +0.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.768"
$ws.Range("E27").Value = "  +3.21%  "

# Row 28
$ws.Range("E28").Value = "  -2.24%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.44"
$ws.Range("E29").Value = "  -0.68%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.939"
$ws.Range("E30").Value = "  -0.44%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07979"
$ws.Range("E31").Value = "  +0.32%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.671"
$ws.Range("E32").Value = "  -0.22%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04485"
$ws.Range("E33").Value = "  -0.89%  "

# Row 34
$ws.Range("E34").Value = "  -0.24%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.608"
$ws.Range("E35").Value = "  -0.19%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.003"
$ws.Range("E36").Value = "  +0.93%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6231"
$ws.Range("E37").Value = "  +0.42%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9436"
$ws.Range("E38").Value = "  +4.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.051"
$ws.Range("E39").Value = "  +3.74%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.425"
$ws.Range("E40").Value = "  +1.73%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9977"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01509"
$ws.Range("E42").Value = "  +1.70%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.571"
$ws.Range("E43").Value = "  +3.57%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.53"
$ws.Range("E44").Value = "  -0.89%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3861"
$ws.Range("E45").Value = "  +0.46%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.906"
$ws.Range("E46").Value = "  +2.77%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1161"
$ws.Range("E47").Value = "  +0.93%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05381"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.891"
$ws.Range("E49").Value = "  +2.29%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.28"
$ws.Range("E50").Value = "  +0.70%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.71"
$ws.Range("E51").Value = "  +1.42%  "
